# Enabling all testcase ECTEST KIT
# Set Runmode (column E) to "Yes" for rows 3 through 29 on the MasterExecutor sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

$ws.Range("E3:E29").Value = "Yes"

# Update the active selection to match the new selection state (E2:E29)
$ws.Range("E2:E29").Select()
